# Re-applies the refreshed cryptocurrency price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.047.72"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "1.676.10"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "'329.59"
$ws.Range("E5").Value = "  +7.09%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").Value = "'47.45"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.3251"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").Value = "'0.07262"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("D12").Value = "'0.9994"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "1.678.17"
$ws.Range("D16").Value = "'6.673"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'0.00001054"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'0.06540"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'0.9988"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'79.08"
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'5.920"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "'12.84"
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").Value = "25.034.83"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "'2.438"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").Value = "'2.399"
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("D27").Value = "'149.28"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "1.866.23"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "'126.12"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").Value = "'1.197"
$ws.Range("D32").Value = "'4.085"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "'5.825"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").Value = "'0.08469"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").Value = "'1.668"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "'5.182"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'0.06108"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02241"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.232"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").Value = "'0.2089"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").Value = "'8.335"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").Value = "'0.9984"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "'0.5999"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").Value = "'13.66"
$ws.Range("E45").Value = "  +8.15%  "
$ws.Range("D46").Value = "'3.833"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").Value = "'0.5750"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "'124.53"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").Value = "'1.968"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "'0.07021"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "'1.188"
$ws.Range("E51").Value = "  +3.13%  "
